$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "65.131.32"
$ws.Range("E2").Value = "  +3.09%  "
$ws.Range("D3").Value = "2.638.36"
$ws.Range("E3").Value = "  +2.88%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.24"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.23"
$ws.Range("E6").Value = "  +5.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("E9").Value = "  +8.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.402"
$ws.Range("E10").Value = "  +5.33%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.154"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.23"
$ws.Range("E13").Value = "  +6.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000187"
$ws.Range("E14").Value = "  +22.06%  "
$ws.Range("D15").Value = "3.114.97"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "65.062.50"
$ws.Range("E16").Value = "  +3.10%  "
$ws.Range("D17").Value = "2.665.20"
$ws.Range("E17").Value = "  +3.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.58"
$ws.Range("E18").Value = "  +3.80%  "
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.37"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("E21").Value = "  +8.30%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.41"
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.56"
$ws.Range("E25").Value = "  +5.53%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.165"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").Value = "0.0₃0951"
$ws.Range("E29").Value = "  +11.22%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +4.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "511.12"
$ws.Range("E32").Value = "  -7.70%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  +8.60%  "
$ws.Range("E35").Value = "  +6.61%  "
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.34"
$ws.Range("E37").Value = "  +5.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.33"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("E39").Value = "  +6.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.30"
$ws.Range("E42").Value = "  +6.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.78"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.09"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.13"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("E47").Value = "  +8.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.648"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0256"
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.50"
$ws.Range("E51").Value = "  +3.15%  "
